$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(0.353672031788087, -0.1145335553579903, 0.1870194611415741, 0.1745300346666848, 0.7152945399284363, 0.1680716276168823, 0.8204078674316406, 0.4750536680221558)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
